$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The script re-ran: the previous "latest" row (row 5, 2025-03-05) no longer
# carries the "NA" marker in column C (it now matches the empty-result rows
# above it), and a brand-new row is appended for 2025-03-06 that carries the
# "NA" marker instead.
$ws.Range("C5").Value = "'"

$ws.Range("A6").Value = "'2025-03-06"
$ws.Range("B6").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C6").Value = "NA"
$ws.Range("D6").Value = 1
